$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 - shifts old rows 13..22 down to 14..23,
# inheriting the row-above style for column A (s=1) automatically.
$ws.Rows.Item(13).Insert()

# The old row 22 content ("Bibliografia principal: ...") moved to row 23 and is
# no longer needed in the new layout (it is dropped entirely) - delete it.
$ws.Rows.Item(23).Delete()

# --- Update the surviving rows 13-21 with their new content ---

# Row 13: "Programa resumido:" / "Semestral"
$ws.Range("A13").Value = "Programa resumido:"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"
$ws.Rows.Item(13).RowHeight = 60

# Row 14: "Short syllabus:" / English short syllabus text
$ws.Range("A14").Value = "Short syllabus:"
$ws.Range("B14").Value = "Introduction to technical drawing software. Graphic primitives and handling operations. 2D drawings. 3D drawings. Import and export of data. Special applications. Kitting. Motion preview and interference."
$ws.Range("C14").Value = "Introduction to technical drawing software. Graphic primitives and handling operations. 2D drawings. 3D drawings. Import and export of data. Special applications. Kitting. Motion preview and interference."
$ws.Rows.Item(14).RowHeight = 60

# Row 15: "Programa:" / "01/01/2018"
$ws.Range("A15").Value = "Programa:"
$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"
$ws.Rows.Item(15).RowHeight = 120

# Row 16: "Syllabus:" / English syllabus text (unchanged content, row renumbered)
$ws.Range("A16").Value = "Syllabus:"
$ws.Range("B16").Value = "Introduction to graphics hardware and peripherals; Generic structure of technical drawing software; Primitive graphics and graphics handling operations; Execution of drawings with technical drawing software: 2D drawings, 3D drawings, views from 3D models, data import and export, use of graphics Libraries, sketches and parametric modeling, special applications, kitting, visualization of movements and interference."
$ws.Range("C16").Value = "Introduction to graphics hardware and peripherals; Generic structure of technical drawing software; Primitive graphics and graphics handling operations; Execution of drawings with technical drawing software: 2D drawings, 3D drawings, views from 3D models, data import and export, use of graphics Libraries, sketches and parametric modeling, special applications, kitting, visualization of movements and interference."
$ws.Rows.Item(16).RowHeight = 120

# Row 17: "Avaliação:" only, no B/C content
$ws.Range("A17").Value = "Avaliação:"
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Rows.Item(17).RowHeight = 15

# Row 18: "Método:" / professor id
$ws.Range("A18").Value = "Método:"
$ws.Range("B18").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Range("C18").Value = "8767640 - Eduardo Ferro dos Santos"
$ws.Rows.Item(18).RowHeight = 60

# Row 19: "Critério:" / "Aulas expositivas e práticas."
$ws.Range("A19").Value = "Critério:"
$ws.Range("B19").Value = "Aulas expositivas e práticas."
$ws.Range("C19").Value = "Aulas expositivas e práticas."
$ws.Rows.Item(19).RowHeight = 60

# Row 20: "Norma de recuperação:" / evaluation criteria text
$ws.Range("A20").Value = "Norma de recuperação:"
$ws.Range("B20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Range("C20").Value = "Exercícios de aprendizado e exercícios de avaliação farão parte da composição de notas individuais (NI), com aplicação de trabalhos práticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2"
$ws.Rows.Item(20).RowHeight = 60

# Row 21: "Bibliografia:" / recovery exam text
$ws.Range("A21").Value = "Bibliografia:"
$ws.Range("B21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Range("C21").Value = "A recuperação deverá consistir de uma prova englobando a matéria toda do semestre. - A média final (pós-recuperação) deverá ser composta por uma média simples entre a nota do semestre (nota final) e a da prova de recuperação."
$ws.Rows.Item(21).RowHeight = 120
